$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D updates to stay as text (avoid Excel auto-numeric coercion)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.174.05"
$ws.Range("E2").Value = "  -2.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.181.65"
$ws.Range("E3").Value = "  -7.35%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.79"
$ws.Range("E5").Value = "  -3.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.74"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.181.02"
$ws.Range("E9").Value = "  -7.33%  "

# Row 10
$ws.Range("E10").Value = "  -5.08%  "

# Row 11
$ws.Range("E11").Value = "  -4.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("E12").Value = "  -2.72%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.735.65"
$ws.Range("E13").Value = "  -7.28%  "

# Row 14
$ws.Range("E14").Value = "  +0.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.63"
$ws.Range("E15").Value = "  -3.12%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.172.43"
$ws.Range("E16").Value = "  -2.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -4.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.177.90"
$ws.Range("E18").Value = "  -7.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.68"
$ws.Range("E19").Value = "  -3.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.10"
$ws.Range("E20").Value = "  -4.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.81"
$ws.Range("E21").Value = "  -3.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.18"
$ws.Range("E22").Value = "  -4.75%  "

# Row 23
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.50"
$ws.Range("E24").Value = "  -3.51%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.504"
$ws.Range("E25").Value = "  -4.73%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  -1.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.48"
$ws.Range("E27").Value = "  -1.48%  "

# Row 28
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.71"
$ws.Range("E30").Value = "  +0.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("E32").Value = "  -3.90%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.15"
$ws.Range("E33").Value = "  -5.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.64"
$ws.Range("E34").Value = "  -4.42%  "

# Row 35
$ws.Range("E35").Value = "  -5.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.06"
$ws.Range("E36").Value = "  -1.94%  "

# Row 37
$ws.Range("E37").Value = "  -5.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.803"
$ws.Range("E38").Value = "  -8.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.99"
$ws.Range("E39").Value = "  -8.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.51"
$ws.Range("E40").Value = "  -3.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("E41").Value = "  -3.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.672.08"
$ws.Range("E42").Value = "  -3.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.15"
$ws.Range("E43").Value = "  -5.96%  "

# Row 44
$ws.Range("E44").Value = "  -6.39%  "

# Row 45
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "328.41"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0651"
$ws.Range("E46").Value = "  -3.57%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.64"
$ws.Range("E47").Value = "  -3.77%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.72"
$ws.Range("E48").Value = "  -2.93%  "

# Row 49
$ws.Range("E49").Value = "  -5.87%  "

# Row 50
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$ws.Range("E51").Value = "  -0.05%  "
